# Update the "Förändrad" date column (C2:C5) from 2023-09-01 to 2023-09-05
# (Excel serial date 45170 -> 45174) as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = Get-Date -Year 2023 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
